$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendrier")

# Row 50 becomes the "Gentlemen de Nommay" entry with its new date,
# row 51 becomes what used to be in row 50 ("5e VTT MS Automobile Rixheim").
$ws.Range("A50").Value = "Dim 20 Septembre"
$ws.Range("B50").Value = "Gentlemen de Nommay (épreuve FFC ouverte aux FSGT)*"
$ws.Range("C50").Value = "CCI Nommay"
$ws.Range("D50").Value = "Route"
$ws.Range("E50").Value = "nommay"

$ws.Range("A51").Value = "Sam 26 Septembre"
$ws.Range("B51").Value = "5e VTT MS Automobile Rixheim  "
$ws.Range("C51").Value = "SSOL Habsheim"
$ws.Range("D51").Value = "VTT"
$ws.Range("E51").Value = "wittenheim"

# Row 49 picks up an explicit height in the saved file.
$ws.Rows.Item(49).RowHeight = 15.75

# Update the view so the selection/scroll position reflects the edited area.
$ws.Activate()
$ws.Range("B50").Select()
$excel.ActiveWindow.ScrollRow = 37
